$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 0.1146748099463193
$ws.Range("F2").Value = 0.0295347868110849
$ws.Range("G2").Value = 0.03289813169417464
$ws.Range("H2").Value = 0.02537661497161529
$ws.Range("I2").Value = 0.08208854512755528
$ws.Range("J2").Value = 0.03704535753565166
$ws.Range("K2").Value = 0.01502104979070028
$ws.Range("L2").Value = 0.002420929129633662
$ws.Range("M2").Value = 0.006944994627913285
$ws.Range("N2").Value = 0.0296048454817489
$ws.Range("O2").Value = 0.06055149305777251
$ws.Range("P2").Value = 0.08276297725771638
$ws.Range("Q2").Value = 0.08431186057955035
$ws.Range("R2").Value = 0.03991023931397532
$ws.Range("S2").Value = 0.0708248254974881
$ws.Range("T2").Value = 0.02621621457098029
$ws.Range("U2").Value = 0.03744072721821363
$ws.Range("V2").Value = 0.004078536025236336
$ws.Range("W2").Value = 0.08448710618906395
$ws.Range("X2").Value = 0.07956804781918368
$ws.Range("Y2").Value = 0.05731131789505488
$ws.Range("Z2").Value = 0.04015333152809897
$ws.Range("AA2").Value = 0.01496822328333368
$ws.Range("AB2").Value = 0.05647984459425399
$ws.Range("AC2").Value = 0.1381283694495126
$ws.Range("E3").Value = 0.1268139114834922
$ws.Range("F3").Value = 0.009718147174482005
$ws.Range("G3").Value = 0.03938869268419044
$ws.Range("H3").Value = 0.06126745479562808
$ws.Range("I3").Value = 0.08768972747948021
$ws.Range("J3").Value = 0.09499911463469989
$ws.Range("K3").Value = 0.004608434468340739
$ws.Range("L3").Value = 0.01074971126812352
$ws.Range("M3").Value = 0.005345391945618391
$ws.Range("N3").Value = 0.06323990843376619
$ws.Range("O3").Value = 0.005578906010639163
$ws.Range("P3").Value = 0.006210286404290825
$ws.Range("Q3").Value = 0.03626880082084739
$ws.Range("R3").Value = 0.02359927074192484
$ws.Range("S3").Value = 0.06124517772936003
$ws.Range("T3").Value = 0.0621063200981407
$ws.Range("U3").Value = 0.03865470924313682
$ws.Range("V3").Value = 0.0239322467402135
$ws.Range("W3").Value = 0.07225394895404734
$ws.Range("X3").Value = 0.04769607201329708
$ws.Range("Y3").Value = 0.0521937845753815
$ws.Range("Z3").Value = 0.03122408375678084
$ws.Range("AA3").Value = 0.06623916001943606
$ws.Range("AB3").Value = 0.09579065000817433
$ws.Range("AC3").Value = 0.1481841094432054
$ws.Range("E4").Value = 0.1189919426803979
$ws.Range("F4").Value = 0.06118993389458165
$ws.Range("G4").Value = 0.05212890100328096
$ws.Range("H4").Value = 0.009700790637828654
$ws.Range("I4").Value = 0.01892618386438418
$ws.Range("J4").Value = 0.08265179527763762
$ws.Range("K4").Value = 0.0005819397504783106
$ws.Range("L4").Value = 0.01017114941099485
$ws.Range("M4").Value = 0.0282839979105465
$ws.Range("N4").Value = 0.04422377464108581
$ws.Range("O4").Value = 0.06979443941160596
$ws.Range("P4").Value = 0.06019284017644116
$ws.Range("Q4").Value = 0.06226606422547118
$ws.Range("R4").Value = 0.09486580292384117
$ws.Range("S4").Value = 0.007940293177930336
$ws.Range("T4").Value = 0.044505554646057
$ws.Range("U4").Value = 0.08253752219042662
$ws.Range("V4").Value = 0.04117543673312082
$ws.Range("W4").Value = 0.08455658265490698
$ws.Range("X4").Value = 0.0456648247987806
$ws.Range("Y4").Value = 0.02197787171292951
$ws.Range("Z4").Value = 0.07135976360792852
$ws.Range("AA4").Value = 0.001953166921525454
$ws.Range("AB4").Value = 0.003351370428216249
$ws.Range("AC4").Value = -0.5920754159267158
$ws.Range("E5").Value = 0.1501250618911742
$ws.Range("F5").Value = 0.009757938496868001
$ws.Range("G5").Value = 0.004510063066974755
$ws.Range("H5").Value = 0.01250540119342636
$ws.Range("I5").Value = 0.1044263660440591
$ws.Range("J5").Value = 0.103276929542559
$ws.Range("K5").Value = 0.01659577499833643
$ws.Range("L5").Value = 0.05477335925303728
$ws.Range("M5").Value = 0.02473031511484418
$ws.Range("N5").Value = 0.01890607544821475
$ws.Range("O5").Value = 0.02649226942971073
$ws.Range("P5").Value = 0.0008975183748449918
$ws.Range("Q5").Value = 0.1035109424948999
$ws.Range("R5").Value = 0.005941908109265788
$ws.Range("S5").Value = 0.05695199691968415
$ws.Range("T5").Value = 0.06846983268466077
$ws.Range("U5").Value = 0.0009763096972649681
$ws.Range("V5").Value = 0.01092009331133864
$ws.Range("W5").Value = 0.03927916788728817
$ws.Range("X5").Value = 0.08609964578592932
$ws.Range("Y5").Value = 0.04131964632434937
$ws.Range("Z5").Value = 0.04921737414997199
$ws.Range("AA5").Value = 0.05627083543057855
$ws.Range("AB5").Value = 0.104170236241893
$ws.Range("AC5").Value = -0.6406194812127495
$ws.Range("E6").Value = 0.2369548880601948
$ws.Range("F6").Value = 0.04041909482602562
$ws.Range("G6").Value = 0.002052960138514599
$ws.Range("H6").Value = 0.03703925800450057
$ws.Range("I6").Value = 0.05194068542110938
$ws.Range("J6").Value = 0.007490005596975861
$ws.Range("K6").Value = 0.04321090148972149
$ws.Range("L6").Value = 0.03664211750827184
$ws.Range("M6").Value = 0.009591725646672024
$ws.Range("N6").Value = 0.00878676158512128
$ws.Range("O6").Value = 0.07807376833286905
$ws.Range("P6").Value = 0.07633901297924099
$ws.Range("Q6").Value = 0.03216128645671244
$ws.Range("R6").Value = 0.0700858610989453
$ws.Range("S6").Value = 0.07104503772429314
$ws.Range("T6").Value = 0.01205609268118535
$ws.Range("U6").Value = 0.08001488699770321
$ws.Range("V6").Value = 0.002578267069584806
$ws.Range("W6").Value = 0.03914991735661853
$ws.Range("X6").Value = 0.06124112365394781
$ws.Range("Y6").Value = 0.07994213215002438
$ws.Range("Z6").Value = 0.06898772316526368
$ws.Range("AA6").Value = 0.05895923088388808
$ws.Range("AB6").Value = 0.03219214923281048
$ws.Range("AC6").Value = 0.1794184707456276
$ws.Range("E7").Value = 0.3565388703949427
$ws.Range("F7").Value = 0.02843619308056724
$ws.Range("G7").Value = 0.1037908952097061
$ws.Range("H7").Value = 0.02732235841592715
$ws.Range("I7").Value = 0.009871803513810692
$ws.Range("J7").Value = 0.02206143498321235
$ws.Range("K7").Value = 0.005109895093482352
$ws.Range("L7").Value = 0.005908658559813168
$ws.Range("M7").Value = 0.00154542539799433
$ws.Range("N7").Value = 0.08505375150436797
$ws.Range("O7").Value = 0.008274046413902831
$ws.Range("P7").Value = 0.02981173740139523
$ws.Range("Q7").Value = 0.06586020738331151
$ws.Range("R7").Value = 0.04609648496973202
$ws.Range("S7").Value = 0.01363207946153071
$ws.Range("T7").Value = 0.005844784687894643
$ws.Range("U7").Value = 0.0002770523162383141
$ws.Range("V7").Value = 0.1092000451887199
$ws.Range("W7").Value = 0.02629252706091982
$ws.Range("X7").Value = 0.1034138639462326
$ws.Range("Y7").Value = 0.09699632087481502
$ws.Range("Z7").Value = 0.08606448853189418
$ws.Range("AA7").Value = 0.05495145674211192
$ws.Range("AB7").Value = 0.06418448926241992
$ws.Range("AC7").Value = 0.3173572464737374
$ws.Range("E8").Value = 0.2466472898311597
$ws.Range("F8").Value = 0.07153631729688154
$ws.Range("G8").Value = 0.02232243119214923
$ws.Range("H8").Value = 0.007263390570059422
$ws.Range("I8").Value = 0.08195564939590774
$ws.Range("J8").Value = 0.06502298962501826
$ws.Range("K8").Value = 0.07469289706335751
$ws.Range("L8").Value = 0.007682784119696247
$ws.Range("M8").Value = 0.04951509607102683
$ws.Range("N8").Value = 0.006597381908475871
$ws.Range("O8").Value = 0.02178826845765678
$ws.Range("P8").Value = 0.06298947276834217
$ws.Range("Q8").Value = 0.01402270036473967
$ws.Range("R8").Value = 0.08915067054468442
$ws.Range("S8").Value = 0.07083176368718351
$ws.Range("T8").Value = 0.007340018970964299
$ws.Range("U8").Value = 0.003123960597102813
$ws.Range("V8").Value = 0.02765686258212393
$ws.Range("W8").Value = 0.08062874819689203
$ws.Range("X8").Value = 0.08446813044348098
$ws.Range("Y8").Value = 0.03198668971629635
$ws.Range("Z8").Value = 0.06285822824453663
$ws.Range("AA8").Value = 0.01030526989135377
$ws.Range("AB8").Value = 0.0462602782920699
$ws.Range("AC8").Value = 0.1566954572801168
$ws.Range("E9").Value = 0.302300845351491
$ws.Range("F9").Value = 0.06151654422431328
$ws.Range("G9").Value = 0.0090109601483423
$ws.Range("H9").Value = 0.03947456837322789
$ws.Range("I9").Value = 0.05062906055923213
$ws.Range("J9").Value = 0.06426848080679341
$ws.Range("K9").Value = 0.1035107339111488
$ws.Range("L9").Value = 0.02772531961778522
$ws.Range("M9").Value = 0.01504239255546109
$ws.Range("N9").Value = 0.00814368943508304
$ws.Range("O9").Value = 0.008359462989723253
$ws.Range("P9").Value = 0.02537434183146755
$ws.Range("Q9").Value = 0.005239601235527331
$ws.Range("R9").Value = 0.03523788409109425
$ws.Range("S9").Value = 0.01261048124134066
$ws.Range("T9").Value = 0.06539488809225126
$ws.Range("U9").Value = 0.0287540882188061
$ws.Range("V9").Value = 0.1287164615879603
$ws.Range("W9").Value = 0.03613659093330603
$ws.Range("X9").Value = 0.01849030984448263
$ws.Range("Y9").Value = 0.1194608225705319
$ws.Range("Z9").Value = 0.05239729363774615
$ws.Range("AA9").Value = 0.007279017210730982
$ws.Range("AB9").Value = 0.07722700688364427
$ws.Range("AC9").Value = 0.2129221450538215
$ws.Range("E10").Value = 0.161786407818635
$ws.Range("F10").Value = 0.08618969168979279
$ws.Range("G10").Value = 0.006243117505752791
$ws.Range("H10").Value = 0.003028070365379065
$ws.Range("I10").Value = 0.0530807005807614
$ws.Range("J10").Value = 0.02696822646664411
$ws.Range("K10").Value = 0.02259461188207384
$ws.Range("L10").Value = 0.003000304423096324
$ws.Range("M10").Value = 0.03130208465394879
$ws.Range("N10").Value = 0.01767264876959547
$ws.Range("O10").Value = 0.07177440285593706
$ws.Range("P10").Value = 0.05552530894900969
$ws.Range("Q10").Value = 0.07311636335210199
$ws.Range("R10").Value = 0.08658977804614991
$ws.Range("S10").Value = 0.05751235530766884
$ws.Range("T10").Value = 0.003087785884997647
$ws.Range("U10").Value = 0.06716769782460116
$ws.Range("V10").Value = 0.0001107432973739428
$ws.Range("W10").Value = 0.06075520584499211
$ws.Range("X10").Value = 0.02423273064502503
$ws.Range("Y10").Value = 0.06863652746741514
$ws.Range("Z10").Value = 0.07248428361864996
$ws.Range("AA10").Value = 0.02263494620603885
$ws.Range("AB10").Value = 0.08629241436299417
$ws.Range("AC10").Value = -0.08114545527317552
$ws.Range("E11").Value = 0.2251703447273393
$ws.Range("F11").Value = 0.09067163847904856
$ws.Range("G11").Value = 0.08050312738725211
$ws.Range("H11").Value = 0.08745204572221939
$ws.Range("I11").Value = 0.00244577079082114
$ws.Range("J11").Value = 0.01765656188882709
$ws.Range("K11").Value = 0.009573655691366962
$ws.Range("L11").Value = 0.09625092099278829
$ws.Range("M11").Value = 0.02015604072834226
$ws.Range("N11").Value = 0.004965503078515965
$ws.Range("O11").Value = 0.004754830782929697
$ws.Range("P11").Value = 0.02008894745581048
$ws.Range("Q11").Value = 0.01168240594717427
$ws.Range("R11").Value = 0.03125969063016845
$ws.Range("S11").Value = 0.07561038683790382
$ws.Range("T11").Value = 0.003836808249559233
$ws.Range("U11").Value = 0.01474591112486043
$ws.Range("V11").Value = 0.09277163442429391
$ws.Range("W11").Value = 0.09804092106167112
$ws.Range("X11").Value = 0.008635479137476816
$ws.Range("Y11").Value = 0.001999834920153905
$ws.Range("Z11").Value = 0.03110194681201118
$ws.Range("AA11").Value = 0.09740114693770735
$ws.Range("AB11").Value = 0.09839479091909763
$ws.Range("AC11").Value = -0.09564825617167652
$ws.Range("E12").Value = 0.1686554608644568
$ws.Range("F12").Value = 0.02435452114359792
$ws.Range("G12").Value = 0.06205927740683664
$ws.Range("H12").Value = 0.01798540610394017
$ws.Range("I12").Value = 0.04284115319791984
$ws.Range("J12").Value = 0.07773661058567646
$ws.Range("K12").Value = 0.001902352719363398
$ws.Range("L12").Value = 0.00193122378555923
$ws.Range("M12").Value = 0.04891494399236151
$ws.Range("N12").Value = 0.003449520493737525
$ws.Range("O12").Value = 0.04120757834327525
$ws.Range("P12").Value = 0.0228108329321995
$ws.Range("Q12").Value = 0.03703623403692852
$ws.Range("R12").Value = 0.09750971302618386
$ws.Range("S12").Value = 0.1084674852007746
$ws.Range("T12").Value = 0.05037006796126385
$ws.Range("U12").Value = 0.1024707472244427
$ws.Range("V12").Value = 0.00338248835052688
$ws.Range("W12").Value = 0.03948396992670766
$ws.Range("X12").Value = 0.06794680207204634
$ws.Range("Y12").Value = 0.01425273605236203
$ws.Range("Z12").Value = 0.05975416009820667
$ws.Range("AA12").Value = 0.05890621316241831
$ws.Range("AB12").Value = 0.01522596218367109
$ws.Range("AC12").Value = 0.2077853665129633
$ws.Range("E13").Value = 0.188559556733204
$ws.Range("F13").Value = 0.09370027855037952
$ws.Range("G13").Value = 0.02903676402788945
$ws.Range("H13").Value = 0.08972902498831442
$ws.Range("I13").Value = 0.07417529730361797
$ws.Range("J13").Value = 0.003424945009904985
$ws.Range("K13").Value = 0.076093541879628
$ws.Range("L13").Value = 0.0003438115135298729
$ws.Range("M13").Value = 0.00001867472556442383
$ws.Range("N13").Value = 0.04076125542068747
$ws.Range("O13").Value = 0.03174557094997298
$ws.Range("P13").Value = 0.003676051600156858
$ws.Range("Q13").Value = 0.001570837053382542
$ws.Range("R13").Value = 0.09293166685583484
$ws.Range("S13").Value = 0.08888457574789004
$ws.Range("T13").Value = 0.002730998928709618
$ws.Range("U13").Value = 0.001810890372181248
$ws.Range("V13").Value = 0.05655077201542323
$ws.Range("W13").Value = 0.04809151039710321
$ws.Range("X13").Value = 0.02712163844192515
$ws.Range("Y13").Value = 0.02203123543650335
$ws.Range("Z13").Value = 0.08660570839419905
$ws.Range("AA13").Value = 0.03869153162954932
$ws.Range("AB13").Value = 0.09027341875765249
$ws.Range("AC13").Value = 0.2455796267029775
$ws.Range("E14").Value = 0.1739911657954868
$ws.Range("F14").Value = 0.1100929545377488
$ws.Range("G14").Value = 0.01177421085815007
$ws.Range("H14").Value = 0.04722999997773143
$ws.Range("I14").Value = 0.01982011485553879
$ws.Range("J14").Value = 0.07150006987690012
$ws.Range("K14").Value = 0.007319367606094846
$ws.Range("L14").Value = 0.0444397488315247
$ws.Range("M14").Value = 0.05258278567386988
$ws.Range("N14").Value = 0.0077618381486169
$ws.Range("O14").Value = 0.00488513828709149
$ws.Range("P14").Value = 0.01114394366240993
$ws.Range("Q14").Value = 0.06454991760990832
$ws.Range("R14").Value = 0.09661206900732998
$ws.Range("S14").Value = 0.07956559453324512
$ws.Range("T14").Value = 0.002552314800963439
$ws.Range("U14").Value = 0.1043543419641532
$ws.Range("V14").Value = 0.001142588426602012
$ws.Range("W14").Value = 0.05970994037600773
$ws.Range("X14").Value = 0.030186962289267
$ws.Range("Y14").Value = 0.03624779270410911
$ws.Range("Z14").Value = 0.06903246113646816
$ws.Range("AA14").Value = 0.06356394721215412
$ws.Range("AB14").Value = 0.003931897624114874
$ws.Range("AC14").Value = 0.1678040005007476
$ws.Range("E15").Value = 0.1981024549833096
$ws.Range("F15").Value = 0.03611027900678178
$ws.Range("G15").Value = 0.02134484757297458
$ws.Range("H15").Value = 0.02005593588387129
$ws.Range("I15").Value = 0.03804270049341741
$ws.Range("J15").Value = 0.006532618646350515
$ws.Range("K15").Value = 0.1037797465231561
$ws.Range("L15").Value = 0.07762396299592944
$ws.Range("M15").Value = 0.03032048842040492
$ws.Range("N15").Value = 0.01180468648185169
$ws.Range("O15").Value = 0.09057641196454418
$ws.Range("P15").Value = 0.04438453804828752
$ws.Range("Q15").Value = 0.01759804405609052
$ws.Range("R15").Value = 0.07107264226588055
$ws.Range("S15").Value = 0.02925456677084153
$ws.Range("T15").Value = 0.03300662282706112
$ws.Range("U15").Value = 0.03761924179535918
$ws.Range("V15").Value = 0.001352747292584261
$ws.Range("W15").Value = 0.03031405719088516
$ws.Range("X15").Value = 0.1221483964424458
$ws.Range("Y15").Value = 0.000935714441678471
$ws.Range("Z15").Value = 0.09810378572408933
$ws.Range("AA15").Value = 0.0202636584277941
$ws.Range("AB15").Value = 0.05775430672772044
$ws.Range("AC15").Value = 0.2428813392719517
$ws.Range("E16").Value = 0.136999816417872
$ws.Range("F16").Value = 0.006853998944420954
$ws.Range("G16").Value = 0.01276777329612985
$ws.Range("H16").Value = 0.06887144169278221
$ws.Range("I16").Value = 0.09211892685638365
$ws.Range("J16").Value = 0.02761680827865268
$ws.Range("K16").Value = 0.004512963732302913
$ws.Range("L16").Value = 0.02360021313016719
$ws.Range("M16").Value = 0.09436494582287418
$ws.Range("N16").Value = 0.007438999619058721
$ws.Range("O16").Value = 0.03106463412875451
$ws.Range("P16").Value = 0.01418348273550831
$ws.Range("Q16").Value = 0.07703474157769773
$ws.Range("R16").Value = 0.1048888875164526
$ws.Range("S16").Value = 0.1148282894334114
$ws.Range("T16").Value = 0.004814627337768751
$ws.Range("U16").Value = 0.07620914467672339
$ws.Range("V16").Value = 0.01413765972235738
$ws.Range("W16").Value = 0.01198112226902973
$ws.Range("X16").Value = 0.005838318849208004
$ws.Range("Y16").Value = 0.03952942661285156
$ws.Range("Z16").Value = 0.08907482252152796
$ws.Range("AA16").Value = 0.06782713735741944
$ws.Range("AB16").Value = 0.01044163388851683
$ws.Range("AC16").Value = 0.05413024478009123
$ws.Range("E17").Value = 0.1622657017886014
$ws.Range("F17").Value = 0.08766906272769012
$ws.Range("G17").Value = 0.01265606487480121
$ws.Range("H17").Value = 0.03133544482658715
$ws.Range("I17").Value = 0.09209046425301136
$ws.Range("J17").Value = 0.03420096710877609
$ws.Range("K17").Value = 0.09849704020933349
$ws.Range("L17").Value = 0.03550103781692002
$ws.Range("M17").Value = 0.06477256989380624
$ws.Range("N17").Value = 0.01311073413796272
$ws.Range("O17").Value = 0.0888184364306656
$ws.Range("P17").Value = 0.05488784685494747
$ws.Range("Q17").Value = 0.04418963849110309
$ws.Range("R17").Value = 0.01379887523969291
$ws.Range("S17").Value = 0.08427508730698383
$ws.Range("T17").Value = 0.02761159401231463
$ws.Range("U17").Value = 0.08041577257421778
$ws.Range("V17").Value = 0.02667557834804073
$ws.Range("W17").Value = 0.02600775327108566
$ws.Range("X17").Value = 0.00004205172350988111
$ws.Range("Y17").Value = 0.01416766817730739
$ws.Range("Z17").Value = 0.002051209296846695
$ws.Range("AA17").Value = 0.06423693481015687
$ws.Range("AB17").Value = 0.002988167614239247
$ws.Range("AC17").Value = 0.02303128585570636
$ws.Range("E18").Value = 0.1229465240068445
$ws.Range("F18").Value = 0.03022481258992275
$ws.Range("G18").Value = 0.04652729239498874
$ws.Range("H18").Value = 0.03774601048718679
$ws.Range("I18").Value = 0.01736614273932583
$ws.Range("J18").Value = 0.04179270730765807
$ws.Range("K18").Value = 0.01034703222792602
$ws.Range("L18").Value = 0.05404371727265311
$ws.Range("M18").Value = 0.01972280808046539
$ws.Range("N18").Value = 0.04760102828545271
$ws.Range("O18").Value = 0.02833704591545615
$ws.Range("P18").Value = 0.05780587862578915
$ws.Range("Q18").Value = 0.03320355737627498
$ws.Range("R18").Value = 0.06457042413424256
$ws.Range("S18").Value = 0.07830445644195741
$ws.Range("T18").Value = 0.01134271871206094
$ws.Range("U18").Value = 0.07729882137219152
$ws.Range("V18").Value = 0.02563391695690887
$ws.Range("W18").Value = 0.004921906830169673
$ws.Range("X18").Value = 0.06969964452397655
$ws.Range("Y18").Value = 0.0747488739047712
$ws.Range("Z18").Value = 0.06408026756648681
$ws.Range("AA18").Value = 0.03086824822396466
$ws.Range("AB18").Value = 0.07381268803017009
$ws.Range("AC18").Value = 0.05181027417346977
$ws.Range("E19").Value = 0.1369919316591015
$ws.Range("F19").Value = 0.01520580658118928
$ws.Range("G19").Value = 0.05257441346756829
$ws.Range("H19").Value = 0.05165829133511202
$ws.Range("I19").Value = 0.0602271365219831
$ws.Range("J19").Value = 0.04312639062140286
$ws.Range("K19").Value = 0.08158658466050846
$ws.Range("L19").Value = 0.09122643428273763
$ws.Range("M19").Value = 0.03463079396281891
$ws.Range("N19").Value = 0.00570800562819863
$ws.Range("O19").Value = 0.03740815316165624
$ws.Range("P19").Value = 0.08516260082236038
$ws.Range("Q19").Value = 0.01774196917404158
$ws.Range("R19").Value = 0.08523749467860783
$ws.Range("S19").Value = 0.008491715669207187
$ws.Range("T19").Value = 0.05799531558991804
$ws.Range("U19").Value = 0.09058375374345413
$ws.Range("V19").Value = 0.07953237381846524
$ws.Range("W19").Value = 0.001543652141434722
$ws.Range("X19").Value = 0.005553059141933011
$ws.Range("Y19").Value = 0.0229376711500455
$ws.Range("Z19").Value = 0.02322107228809974
$ws.Range("AA19").Value = 0.01988379720273776
$ws.Range("AB19").Value = 0.02876351435651956
$ws.Range("AC19").Value = 0.1426487826520407
$ws.Range("E20").Value = 0.1574065361859451
$ws.Range("F20").Value = 0.05103375427726918
$ws.Range("G20").Value = 0.03621858256461116
$ws.Range("H20").Value = 0.1085297796841402
$ws.Range("I20").Value = 0.008020943596369337
$ws.Range("J20").Value = 0.006530805331893473
$ws.Range("K20").Value = 0.006954940537084406
$ws.Range("L20").Value = 0.02962173580982032
$ws.Range("M20").Value = 0.03375519240435478
$ws.Range("N20").Value = 0.04829023755735993
$ws.Range("O20").Value = 0.07067958205406742
$ws.Range("P20").Value = 0.01037966146423374
$ws.Range("Q20").Value = 0.04291535756821583
$ws.Range("R20").Value = 0.08342436978103138
$ws.Range("S20").Value = 0.008090512086322171
$ws.Range("T20").Value = 0.01948830762447873
$ws.Range("U20").Value = 0.111251464798295
$ws.Range("V20").Value = 0.0244972579787592
$ws.Range("W20").Value = 0.0355467424543421
$ws.Range("X20").Value = 0.08648263776175918
$ws.Range("Y20").Value = 0.02140790876456006
$ws.Range("Z20").Value = 0.04180435963143015
$ws.Range("AA20").Value = 0.1038540351203364
$ws.Range("AB20").Value = 0.01122183114926587
$ws.Range("AC20").Value = 0.1309168858009779
$ws.Range("E21").Value = 0.1733050686683223
$ws.Range("F21").Value = 0.1040709204445876
$ws.Range("G21").Value = 0.01202874380128892
$ws.Range("H21").Value = 0.01546653482907819
$ws.Range("I21").Value = 0.02962469611171085
$ws.Range("J21").Value = 0.08844008788503721
$ws.Range("K21").Value = 0.1153689260004395
$ws.Range("L21").Value = 0.008840662608493243
$ws.Range("M21").Value = 0.04962042318658743
$ws.Range("N21").Value = 0.004645693661523366
$ws.Range("O21").Value = 0.01888488853692425
$ws.Range("P21").Value = 0.04546496830550752
$ws.Range("Q21").Value = 0.001693219591531406
$ws.Range("R21").Value = 0.02385846387167331
$ws.Range("S21").Value = 0.08111482413267204
$ws.Range("T21").Value = 0.03011639327029646
$ws.Range("U21").Value = 0.009428319763379921
$ws.Range("V21").Value = 0.1201617138921037
$ws.Range("W21").Value = 0.06993795084499986
$ws.Range("X21").Value = 0.03955421126353247
$ws.Range("Y21").Value = 0.03495077523610816
$ws.Range("Z21").Value = 0.03143794850593563
$ws.Range("AA21").Value = 0.05333116482529551
$ws.Range("AB21").Value = 0.01195846943129341
$ws.Range("AC21").Value = 0.1229992669726499
$ws.Range("E22").Value = 0.1721954137966594
$ws.Range("F22").Value = 0.04666940809607296
$ws.Range("G22").Value = 0.009086679095127293
$ws.Range("H22").Value = 0.07585641696592287
$ws.Range("I22").Value = 0.04762588228639419
$ws.Range("J22").Value = 0.03709597360383476
$ws.Range("K22").Value = 0.01610911666859032
$ws.Range("L22").Value = 0.0138433350052759
$ws.Range("M22").Value = 0.1000246441598572
$ws.Range("N22").Value = 0.007073017500259633
$ws.Range("O22").Value = 0.01954070428766681
$ws.Range("P22").Value = 0.01348759577135906
$ws.Range("Q22").Value = 0.02051169927081333
$ws.Range("R22").Value = 0.0724661851514069
$ws.Range("S22").Value = 0.0528416966124293
$ws.Range("T22").Value = 0.02848429270526195
$ws.Range("U22").Value = 0.07758231291352623
$ws.Range("V22").Value = 0.01935972280896095
$ws.Range("W22").Value = 0.071418203167656
$ws.Range("X22").Value = 0.08301734555597146
$ws.Range("Y22").Value = 0.03358532118878982
$ws.Range("Z22").Value = 0.09513691853011298
$ws.Range("AA22").Value = 0.05050716622525027
$ws.Range("AB22").Value = 0.008676362429459876
$ws.Range("AC22").Value = 0.05975216194118629
$ws.Range("E23").Value = 0.18402115395109
$ws.Range("F23").Value = 0.08838879784282458
$ws.Range("G23").Value = 0.01854254932030156
$ws.Range("H23").Value = 0.04423985722342036
$ws.Range("I23").Value = 0.05822249219222786
$ws.Range("J23").Value = 0.010879021014854
$ws.Range("K23").Value = 0.03655017615341621
$ws.Range("L23").Value = 0.03072559602237531
$ws.Range("M23").Value = 0.01272183079099692
$ws.Range("N23").Value = 0.01769406691328479
$ws.Range("O23").Value = 0.01419572250386735
$ws.Range("P23").Value = 0.03279432712792568
$ws.Range("Q23").Value = 0.0155625401236394
$ws.Range("R23").Value = 0.02919513804998235
$ws.Range("S23").Value = 0.01503872092474187
$ws.Range("T23").Value = 0.05626002390013125
$ws.Range("U23").Value = 0.078434295664915
$ws.Range("V23").Value = 0.0755555428943912
$ws.Range("W23").Value = 0.05860014006027695
$ws.Range("X23").Value = 0.06470166371572407
$ws.Range("Y23").Value = 0.08907464276036625
$ws.Range("Z23").Value = 0.08629111141786969
$ws.Range("AA23").Value = 0.06140597580003988
$ws.Range("AB23").Value = 0.004925767582427437
$ws.Range("AC23").Value = 0.08451691517183725
$ws.Range("E24").Value = 0.1271196225536832
$ws.Range("F24").Value = 0.08757588238961329
$ws.Range("G24").Value = 0.04393964342166121
$ws.Range("H24").Value = 0.03916740348283879
$ws.Range("I24").Value = 0.008303395120030255
$ws.Range("J24").Value = 0.008087065190474498
$ws.Range("K24").Value = 0.03344685922028771
$ws.Range("L24").Value = 0.05716838767963366
$ws.Range("M24").Value = 0.07298217135115209
$ws.Range("N24").Value = 0.01359464656129023
$ws.Range("O24").Value = 0.03845923969252493
$ws.Range("P24").Value = 0.05700357290224868
$ws.Range("Q24").Value = 0.01614536497886022
$ws.Range("R24").Value = 0.06174497545900642
$ws.Range("S24").Value = 0.0424787989007642
$ws.Range("T24").Value = 0.03453636425940534
$ws.Range("U24").Value = 0.06570573678548332
$ws.Range("V24").Value = 0.001943427457105809
$ws.Range("W24").Value = 0.08922631994253172
$ws.Range("X24").Value = 0.03501782648883259
$ws.Range("Y24").Value = 0.02376467551201793
$ws.Range("Z24").Value = 0.07955809886896863
$ws.Range("AA24").Value = 0.08392721087306364
$ws.Range("AB24").Value = 0.00622293346220492
$ws.Range("AC24").Value = -0.1766180650056217
$ws.Range("E25").Value = 0.1452320568075288
$ws.Range("F25").Value = 0.08625321973011982
$ws.Range("G25").Value = 0.07087662153954116
$ws.Range("H25").Value = 0.001164754226189505
$ws.Range("I25").Value = 0.05009379527604529
$ws.Range("J25").Value = 0.05130398452193009
$ws.Range("K25").Value = 0.01323145925632087
$ws.Range("L25").Value = 0.02133317757094453
$ws.Range("M25").Value = 0.07199489020900923
$ws.Range("N25").Value = 0.02559551051769914
$ws.Range("O25").Value = 0.001138311028389126
$ws.Range("P25").Value = 0.0301723349245401
$ws.Range("Q25").Value = 0.1116651484649404
$ws.Range("R25").Value = 0.07478888258992319
$ws.Range("S25").Value = 0.002227140059501185
$ws.Range("T25").Value = 0.01494857126762505
$ws.Range("U25").Value = 0.06116883410780519
$ws.Range("V25").Value = 0.05529600051433787
$ws.Range("W25").Value = 0.01007910811772871
$ws.Range("X25").Value = 0.006429484334478026
$ws.Range("Y25").Value = 0.0758737774050339
$ws.Range("Z25").Value = 0.07023662919701411
$ws.Range("AA25").Value = 0.03275852537931429
$ws.Range("AB25").Value = 0.06136983976156919
$ws.Range("AC25").Value = -0.1523574155980028
$ws.Range("E26").Value = 0.1080448830276254
$ws.Range("F26").Value = 0.06731776969204754
$ws.Range("G26").Value = 0.07972057556193624
$ws.Range("H26").Value = 0.04319429809983217
$ws.Range("I26").Value = 0.0473018014081226
$ws.Range("J26").Value = 0.05169875360982497
$ws.Range("K26").Value = 0.01790039860717675
$ws.Range("L26").Value = 0.02676692727200177
$ws.Range("M26").Value = 0.08358528787295025
$ws.Range("N26").Value = 0.02681928059779367
$ws.Range("O26").Value = 0.02393387360861241
$ws.Range("P26").Value = 0.05774561772611354
$ws.Range("Q26").Value = 0.08561515691590814
$ws.Range("R26").Value = 0.01322514347330648
$ws.Range("S26").Value = 0.01239805029141942
$ws.Range("T26").Value = 0.006241271119404993
$ws.Range("U26").Value = 0.0390429119471419
$ws.Range("V26").Value = 0.003020958774167717
$ws.Range("W26").Value = 0.08416667046303522
$ws.Range("X26").Value = 0.09517878760270281
$ws.Range("Y26").Value = 0.0009869223789110012
$ws.Range("Z26").Value = 0.08326145303182467
$ws.Range("AA26").Value = 0.0001980136037845736
$ws.Range("AB26").Value = 0.05068007634198141
$ws.Range("AC26").Value = 0.1897404442158593
$ws.Range("E27").Value = 0.1148914013378153
$ws.Range("F27").Value = 0.08439379932595822
$ws.Range("G27").Value = 0.07367723853333913
$ws.Range("H27").Value = 0.06866768112063243
$ws.Range("I27").Value = 0.02391886317048711
$ws.Range("J27").Value = 0.01336561012049093
$ws.Range("K27").Value = 0.03684920237567637
$ws.Range("L27").Value = 0.007789940229888511
$ws.Range("M27").Value = 0.09306293853703841
$ws.Range("N27").Value = 0.002169455004016108
$ws.Range("O27").Value = 0.0105933779796057
$ws.Range("P27").Value = 0.03047132239469818
$ws.Range("Q27").Value = 0.07201884023287156
$ws.Range("R27").Value = 0.03981240493099625
$ws.Range("S27").Value = 0.001274788322287553
$ws.Range("T27").Value = 0.002037837707713206
$ws.Range("U27").Value = 0.03383242520272158
$ws.Range("V27").Value = 0.06493761116217205
$ws.Range("W27").Value = 0.07756661059560406
$ws.Range("X27").Value = 0.0434522377351834
$ws.Range("Y27").Value = 0.08505170485976324
$ws.Range("Z27").Value = 0.0763730216336934
$ws.Range("AA27").Value = 0.007071010310314748
$ws.Range("AB27").Value = 0.05161207851484785
$ws.Range("AC27").Value = 0.2098866037425231
$ws.Range("E28").Value = 0.1138570349336037
$ws.Range("F28").Value = 0.03950903238417761
$ws.Range("G28").Value = 0.05548481904426953
$ws.Range("H28").Value = 0.0742375461996268
$ws.Range("I28").Value = 0.0007889096306436646
$ws.Range("J28").Value = 0.07185855198213975
$ws.Range("K28").Value = 0.00581767035742223
$ws.Range("L28").Value = 0.008108802240525202
$ws.Range("M28").Value = 0.07379528294616701
$ws.Range("N28").Value = 0.00956183895398261
$ws.Range("O28").Value = 0.05102487163382042
$ws.Range("P28").Value = 0.07695602674440873
$ws.Range("Q28").Value = 0.07487911738811016
$ws.Range("R28").Value = 0.02156660535022785
$ws.Range("S28").Value = 0.04997819402561147
$ws.Range("T28").Value = 0.04046289853722759
$ws.Range("U28").Value = 0.07124688419764989
$ws.Range("V28").Value = 0.0009409486107140393
$ws.Range("W28").Value = 0.08498968068035946
$ws.Range("X28").Value = 0.04488923005828022
$ws.Range("Y28").Value = 0.002463591410747997
$ws.Range("Z28").Value = 0.04806538759371646
$ws.Range("AA28").Value = 0.05242173923372292
$ws.Range("AB28").Value = 0.04095237079644835
$ws.Range("AC28").Value = -0.01039839157509783
$ws.Range("E29").Value = 0.1284227210582433
$ws.Range("F29").Value = 0.07162621645916044
$ws.Range("G29").Value = 0.05345028330596475
$ws.Range("H29").Value = 0.07974273425330437
$ws.Range("I29").Value = 0.01255602195526663
$ws.Range("J29").Value = 0.005245538468829529
$ws.Range("K29").Value = 0.01456769453202144
$ws.Range("L29").Value = 0.002158045012698002
$ws.Range("M29").Value = 0.09760291623649349
$ws.Range("N29").Value = 0.01705716938615304
$ws.Range("O29").Value = 0.01951676871079585
$ws.Range("P29").Value = 0.02058702082959307
$ws.Range("Q29").Value = 0.02205380343888732
$ws.Range("R29").Value = 0.08288606027863675
$ws.Range("S29").Value = 0.05554632292395451
$ws.Range("T29").Value = 0.01919308739209181
$ws.Range("U29").Value = 0.00873883095601777
$ws.Range("V29").Value = 0.06799953491363293
$ws.Range("W29").Value = 0.090368828569002
$ws.Range("X29").Value = 0.02254378866763983
$ws.Range("Y29").Value = 0.0929886766653836
$ws.Range("Z29").Value = 0.08651762532060081
$ws.Range("AA29").Value = 0.03170856928886195
$ws.Range("AB29").Value = 0.02534446243501008
$ws.Range("AC29").Value = 0.03725690099284405
